$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")
$ws.Range("C2").Value = 1.403
$ws.Range("D2").Value = 0.1
$ws.Range("F2").Value = "from CEA, costs in USD-2015,"
$ws.Range("F16").Select()
